$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.478.00"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "2.331.75"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "545.72"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "131.43"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "2.329.43"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "5.53"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "2.747.39"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "60.403.30"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D18").Value = "2.332.41"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "314.71"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").Value = "  +5.93%  "
$ws.Range("E29").Value = "  +9.02%  "
$ws.Range("D30").Value = "172.67"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  +9.60%  "
$ws.Range("D35").Value = "0.381"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "321.34"
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("D41").Value = "1.54"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "37.91"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "137.39"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "0.0494"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "0.0214"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "0.0₆0216"
$ws.Range("E50").Value = "  +15.66%  "
$ws.Range("D51").Value = "11.02"
$ws.Range("E51").Value = "  +0.58%  "
